$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 6 new rows after the existing data row (row 16), pushing the
#     footer block (rows 21-22) down to rows 27-28 ---
$ws.Rows.Item(17).Resize(6).Insert()

# Copy formatting (borders, fonts, number formats) from the existing data
# row (16) down into the 6 freshly inserted rows (17-22)
for ($i = 17; $i -le 22; $i++) {
    $ws.Range("B16:J16").Copy($ws.Range("B$i`:J$i"))
}

# --- Fill the new worker's rows (OSWALDO ENRIQUE CARBALLO DE VOZ) ---
$periods = @(2507, 2506, 2505, 2504, 2503, 2502)
$row = 17
foreach ($p in $periods) {
    $ws.Range("B$row").Value = "CC"
    $ws.Range("C$row").Value = "1044917171"
    $ws.Range("D$row").Value = "OSWALDO ENRIQUE CARBALLO DE VOZ"
    $ws.Range("E$row").Value = [string]$p
    $ws.Range("F$row").Value = 36341
    $ws.Range("G$row").Value = 1300000
    $row++
}

# --- Close off the bottom border of the table on the final data row ---
$ws.Range("B22:J22").Borders.Item(9).LineStyle = 1
$ws.Range("B22:J22").Borders.Item(9).Weight = 2
$ws.Range("B22:J22").Borders.Item(9).ColorIndex = 0

# --- Update the summary header values ---
$ws.Range("E11").Value = 220469
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 7

# --- Re-fit column D now that it holds a longer worker name ---
$ws.Columns.Item(4).AutoFit()
